$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCommitment")

# Insert a new column before column I (9th column), shifting existing I:N columns to J:O
$ws.Columns("I:I").Insert()

# Set the new column's header and values
$ws.Range("I1").Value = "Onboarding Completed"
$ws.Range("I2").Value = "Yes"
$ws.Range("I3").Value = "Yes"
$ws.Range("I4").Value = "Yes"
$ws.Range("I5").Value = "Yes"
$ws.Range("I6").Value = "Yes"
$ws.Range("I7").Value = "Yes"

# Update selection to I2 as in the target workbook
$ws.Range("I2").Select()
